$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1485.7333
$ws.Range("I100").Value = 1445.0358
$ws.Range("J100").Value = 2055.5
$ws.Range("K100").Value = 1445.0358
$ws.Range("L100").Value = 2055.5
$ws.Range("M100").Value = -904.0358000000001
$ws.Range("N100").Value = -3137.5
$ws.Range("H103").Value = 290.0625
$ws.Range("I103").Value = 212.71428
$ws.Range("J103").Value = 350.22223
$ws.Range("K103").Value = 638.14284
$ws.Range("L103").Value = 1050.66669
$ws.Range("M103").Value = -52.14283999999998
$ws.Range("N103").Value = -2222.66669
$ws.Range("H107").Value = 777.32355
$ws.Range("I107").Value = 748.5172
$ws.Range("J107").Value = 944.4
$ws.Range("K107").Value = 748.5172
$ws.Range("L107").Value = 944.4
$ws.Range("M107").Value = 1171.4828
$ws.Range("N107").Value = -4784.4
$ws.Range("H113").Value = 15940.091
$ws.Range("I113").Value = 17068.4
$ws.Range("J113").Value = 14999.833
$ws.Range("K113").Value = 17068.4
$ws.Range("L113").Value = 14999.833
$ws.Range("M113").Value = -13814.4
$ws.Range("N113").Value = -21507.833
$ws.Range("H132").Value = 2211.2856
$ws.Range("I132").Value = 2150.641
$ws.Range("J132").Value = 2999.6667
$ws.Range("K132").Value = 6451.923000000001
$ws.Range("L132").Value = 8999.000100000001
$ws.Range("M132").Value = -3921.923000000001
$ws.Range("N132").Value = -14059.0001
$ws.Range("H135").Value = 2605.8
$ws.Range("I135").Value = 2498.25
$ws.Range("K135").Value = 22484.25
$ws.Range("M135").Value = -19949.25
$ws.Range("H137").Value = 4228.625
$ws.Range("I137").Value = 6940
$ws.Range("J137").Value = 2996.182
$ws.Range("K137").Value = 20820
$ws.Range("L137").Value = 8988.545999999998
$ws.Range("M137").Value = -18270
$ws.Range("N137").Value = -14088.546
$ws.Range("H138").Value = 5689.3
$ws.Range("I138").Value = 4056.5
$ws.Range("J138").Value = 5831.2827
$ws.Range("K138").Value = 12169.5
$ws.Range("L138").Value = 17493.8481
$ws.Range("M138").Value = -7029.5
$ws.Range("N138").Value = -27773.8481

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10767.66
$ws.Range("I32").Value = 10171.234
$ws.Range("J32").Value = 39992.5
$ws.Range("K32").Value = 10171.234
$ws.Range("L32").Value = 39992.5
$ws.Range("M32").Value = -9884.234
$ws.Range("N32").Value = -40566.5
$ws.Range("H61").Value = 504811.5
$ws.Range("I61").Value = 3738.25
$ws.Range("K61").Value = 3738.25
$ws.Range("M61").Value = -3526.25
$ws.Range("H74").Value = 54798.09
$ws.Range("I74").Value = 60180.293
$ws.Range("K74").Value = 60180.293
$ws.Range("M74").Value = -59306.293
$ws.Range("H77").Value = 54798.09
$ws.Range("I77").Value = 60180.293
$ws.Range("K77").Value = 300901.465
$ws.Range("M77").Value = -296533.465
$ws.Range("H105").Value = 50366
$ws.Range("J105").Value = 50366
$ws.Range("L105").Value = 50366
$ws.Range("N105").Value = -57354
$ws.Range("H122").Value = 2358.75
$ws.Range("I122").Value = 2358.75
$ws.Range("K122").Value = 7076.25
$ws.Range("M122").Value = -4626.25
$ws.Range("H132").Value = 5532.6
$ws.Range("I132").Value = 5280.722
$ws.Range("J132").Value = 5910.4165
$ws.Range("K132").Value = 15842.166
$ws.Range("L132").Value = 17731.2495
$ws.Range("M132").Value = -13312.166
$ws.Range("N132").Value = -22791.2495
$ws.Range("H136").Value = 504811.5
$ws.Range("I136").Value = 3738.25
$ws.Range("K136").Value = 11214.75
$ws.Range("M136").Value = -8664.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 42652.5
$ws.Range("J28").Value = 42652.5
$ws.Range("L28").Value = 42652.5
$ws.Range("N28").Value = -43240.5
$ws.Range("H86").Value = 2173.182
$ws.Range("I86").Value = 2179.6667
$ws.Range("K86").Value = 2179.6667
$ws.Range("M86").Value = -1056.6667
$ws.Range("H89").Value = 2173.182
$ws.Range("I89").Value = 2179.6667
$ws.Range("K89").Value = 10898.3335
$ws.Range("M89").Value = -5282.333500000001
$ws.Range("H100").Value = 24925.334
$ws.Range("J100").Value = 24925.334
$ws.Range("L100").Value = 24925.334
$ws.Range("N100").Value = -27089.334
$ws.Range("H134").Value = 3073.3333
$ws.Range("I134").Value = 2072.4285
$ws.Range("K134").Value = 6217.2855
$ws.Range("M134").Value = -3682.2855

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2760.3044
$ws.Range("I31").Value = 2613.2273
$ws.Range("J31").Value = 5996
$ws.Range("K31").Value = 2613.2273
$ws.Range("L31").Value = 5996
$ws.Range("M31").Value = -2318.2273
$ws.Range("N31").Value = -6586
$ws.Range("H34").Value = 2760.3044
$ws.Range("I34").Value = 2613.2273
$ws.Range("J34").Value = 5996
$ws.Range("K34").Value = 2613.2273
$ws.Range("L34").Value = 5996
$ws.Range("M34").Value = -2411.2273
$ws.Range("N34").Value = -6400
$ws.Range("H93").Value = 3984.923
$ws.Range("I93").Value = 3984.923
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 3984.923
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -2112.923
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 3380.8667
$ws.Range("I132").Value = 3066.5833
$ws.Range("K132").Value = 9199.749899999999
$ws.Range("M132").Value = -6669.749899999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1249.625
$ws.Range("I5").Value = 703.125
$ws.Range("J5").Value = 2342.625
$ws.Range("K5").Value = 2109.375
$ws.Range("L5").Value = 7027.875
$ws.Range("M5").Value = -1997.375
$ws.Range("N5").Value = -7251.875
$ws.Range("H68").Value = 557741.5600000001
$ws.Range("J68").Value = 1668666.6
$ws.Range("L68").Value = 5005999.800000001
$ws.Range("N68").Value = -5007621.800000001
$ws.Range("H70").Value = 7217.8335
$ws.Range("I70").Value = 5681.4
$ws.Range("K70").Value = 17044.2
$ws.Range("M70").Value = -16729.2
$ws.Range("H71").Value = 557741.5600000001
$ws.Range("J71").Value = 1668666.6
$ws.Range("L71").Value = 15017999.4
$ws.Range("N71").Value = -15026111.4
$ws.Range("H73").Value = 7217.8335
$ws.Range("I73").Value = 5681.4
$ws.Range("K73").Value = 17044.2
$ws.Range("M73").Value = -15952.2
$ws.Range("H129").Value = 9900682
$ws.Range("I129").Value = 12375481
$ws.Range("K129").Value = 37126443
$ws.Range("M129").Value = -37121443
$ws.Range("H135").Value = 1249.625
$ws.Range("I135").Value = 703.125
$ws.Range("J135").Value = 2342.625
$ws.Range("K135").Value = 6328.125
$ws.Range("L135").Value = 21083.625
$ws.Range("M135").Value = -3793.125
$ws.Range("N135").Value = -26153.625
$ws.Range("H137").Value = 3902.1428
$ws.Range("J137").Value = 3966.6667
$ws.Range("L137").Value = 11900.0001
$ws.Range("N137").Value = -22100.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 23948.5
$ws.Range("I52").Value = 23948.5
$ws.Range("K52").Value = 23948.5
$ws.Range("M52").Value = -23689.5
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 2764.739
$ws.Range("I122").Value = 2724.9285
$ws.Range("J122").Value = 2826.6667
$ws.Range("K122").Value = 8174.7855
$ws.Range("L122").Value = 8480.000100000001
$ws.Range("M122").Value = -5724.7855
$ws.Range("N122").Value = -13380.0001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4270.75
$ws.Range("I16").Value = 3369.457
$ws.Range("J16").Value = 10579.8
$ws.Range("K16").Value = 3369.457
$ws.Range("L16").Value = 10579.8
$ws.Range("M16").Value = -3199.457
$ws.Range("N16").Value = -10919.8
$ws.Range("H22").Value = 4481.7646
$ws.Range("I22").Value = 3598.5715
$ws.Range("J22").Value = 5100
$ws.Range("K22").Value = 3598.5715
$ws.Range("L22").Value = 5100
$ws.Range("M22").Value = -3303.5715
$ws.Range("N22").Value = -5690
$ws.Range("H27").Value = 4481.7646
$ws.Range("I27").Value = 3598.5715
$ws.Range("J27").Value = 5100
$ws.Range("K27").Value = 3598.5715
$ws.Range("L27").Value = 5100
$ws.Range("M27").Value = -3491.5715
$ws.Range("N27").Value = -5314
$ws.Range("H46").Value = 3183.5833
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3183.5833
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3183.5833
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -3559.5833
$ws.Range("H82").Value = 5925.04
$ws.Range("I82").Value = 7003.3335
$ws.Range("K82").Value = 7003.3335
$ws.Range("M82").Value = -6642.3335
$ws.Range("H85").Value = 5925.04
$ws.Range("I85").Value = 7003.3335
$ws.Range("K85").Value = 7003.3335
$ws.Range("M85").Value = -5755.3335
$ws.Range("H122").Value = 408126.6
$ws.Range("I122").Value = 674278.8
$ws.Range("J122").Value = 8898.299999999999
$ws.Range("K122").Value = 2022836.4
$ws.Range("L122").Value = 26694.9
$ws.Range("M122").Value = -2020386.4
$ws.Range("N122").Value = -31594.9
$ws.Range("H136").Value = 4952.615
$ws.Range("I136").Value = 4580.364
$ws.Range("K136").Value = 13741.092
$ws.Range("M136").Value = -11191.092

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3305.1765
$ws.Range("I122").Value = 3517.5454
$ws.Range("K122").Value = 10552.6362
$ws.Range("M122").Value = -8102.636200000001
$ws.Range("H136").Value = 2126.2083
$ws.Range("I136").Value = 1910.409
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 5731.227000000001
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -3181.227000000001
$ws.Range("N136").Value = -18600
